# Apply numeric updates to the Golem_Profits leve-profit sheets
# (values recomputed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 41.6
$ws.Range("I6").Value = 21
$ws.Range("J6").Value = 124
$ws.Range("K6").Value = 63
$ws.Range("L6").Value = 372
$ws.Range("M6").Value = 49
$ws.Range("N6").Value = -596
$ws.Range("H8").Value = 10.5
$ws.Range("I8").Value = 1
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 136
$ws.Range("H17").Value = 2626.4
$ws.Range("J17").Value = 2299.7144
$ws.Range("L17").Value = 6899.1432
$ws.Range("N17").Value = -7235.1432
$ws.Range("H28").Value = 1425.8182
$ws.Range("I28").Value = 550.8333
$ws.Range("K28").Value = 550.8333
$ws.Range("M28").Value = -65.83330000000001
$ws.Range("H39").Value = 404.69232
$ws.Range("I39").Value = 442
$ws.Range("J39").Value = 199.5
$ws.Range("K39").Value = 1326
$ws.Range("L39").Value = 598.5
$ws.Range("M39").Value = -1030
$ws.Range("N39").Value = -1190.5
$ws.Range("H42").Value = 453.7143
$ws.Range("I42").Value = 119.75
$ws.Range("J42").Value = 899
$ws.Range("K42").Value = 359.25
$ws.Range("L42").Value = 2697
$ws.Range("M42").Value = -129.25
$ws.Range("N42").Value = -3157
$ws.Range("H43").Value = 6309.1
$ws.Range("I43").Value = 5198
$ws.Range("K43").Value = 5198
$ws.Range("M43").Value = -5129
$ws.Range("H54").Value = 8100
$ws.Range("I54").Value = 8000
$ws.Range("J54").Value = 8400
$ws.Range("K54").Value = 8000
$ws.Range("L54").Value = 8400
$ws.Range("M54").Value = -7514
$ws.Range("N54").Value = -9372
$ws.Range("H76").Value = 3500
$ws.Range("I76").Value = 3500
$ws.Range("K76").Value = 3500
$ws.Range("M76").Value = -3185
$ws.Range("H79").Value = 3500
$ws.Range("I79").Value = 3500
$ws.Range("K79").Value = 3500
$ws.Range("M79").Value = -2408
$ws.Range("H99").Value = 83334150
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 9000
$ws.Range("N99").Value = -11996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 30051
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 55000
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H94").Value = 65628.47
$ws.Range("J94").Value = 966.6667
$ws.Range("L94").Value = 966.6667
$ws.Range("N94").Value = -1868.6667
$ws.Range("H99").Value = 5009.5
$ws.Range("I99").Value = 5009.5
$ws.Range("K99").Value = 5009.5
$ws.Range("M99").Value = -3511.5
$ws.Range("H105").Value = 1984.75
$ws.Range("I105").Value = 1984.75
$ws.Range("K105").Value = 1984.75
$ws.Range("M105").Value = -237.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 40
$ws.Range("M4").Value = 72
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H62").Value = 3166.6667
$ws.Range("I62").Value = 3166.6667
$ws.Range("K62").Value = 3166.6667
$ws.Range("M62").Value = -2542.6667
$ws.Range("H65").Value = 3166.6667
$ws.Range("I65").Value = 3166.6667
$ws.Range("K65").Value = 15833.3335
$ws.Range("M65").Value = -12713.3335
$ws.Range("H134").Value = 1629.3846
$ws.Range("I134").Value = 1471.091
$ws.Range("K134").Value = 4413.272999999999
$ws.Range("M134").Value = -1878.272999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 363.0909
$ws.Range("I13").Value = 120
$ws.Range("J13").Value = 502
$ws.Range("K13").Value = 360
$ws.Range("L13").Value = 1506
$ws.Range("M13").Value = -192
$ws.Range("N13").Value = -1842
$ws.Range("H23").Value = 39.5
$ws.Range("I23").Value = 39.5
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 118.5
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 116.5
$ws.Range("N23").ClearContents()
$ws.Range("H121").Value = 482.5
$ws.Range("I121").Value = 272.4
$ws.Range("J121").Value = 1007.75
$ws.Range("K121").Value = 817.1999999999999
$ws.Range("L121").Value = 3023.25
$ws.Range("M121").Value = 492.8000000000001
$ws.Range("N121").Value = -5643.25
$ws.Range("H131").Value = 2356.875
$ws.Range("I131").Value = 1215
$ws.Range("K131").Value = 3645
$ws.Range("M131").Value = 1395
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 62000
$ws.Range("J39").Value = 62000
$ws.Range("L39").Value = 62000
$ws.Range("N39").Value = -63064
$ws.Range("H70").Value = 29416320
$ws.Range("I70").Value = 3499.5
$ws.Range("J70").Value = 33338028
$ws.Range("K70").Value = 3499.5
$ws.Range("L70").Value = 33338028
$ws.Range("M70").Value = -3229.5
$ws.Range("N70").Value = -33338568
$ws.Range("H73").Value = 29416320
$ws.Range("I73").Value = 3499.5
$ws.Range("J73").Value = 33338028
$ws.Range("K73").Value = 3499.5
$ws.Range("L73").Value = 33338028
$ws.Range("M73").Value = -2563.5
$ws.Range("N73").Value = -33339900
$ws.Range("H80").Value = 4800
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4800
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4800
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -6796
$ws.Range("H83").Value = 4800
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4800
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 24000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -33984
$ws.Range("H113").Value = 972.5
$ws.Range("I113").Value = 972.5
$ws.Range("K113").Value = 972.5
$ws.Range("M113").Value = 1197.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 442.5
$ws.Range("I22").Value = 385
$ws.Range("K22").Value = 385
$ws.Range("M22").Value = -90
$ws.Range("H27").Value = 442.5
$ws.Range("I27").Value = 385
$ws.Range("K27").Value = 385
$ws.Range("M27").Value = -278
$ws.Range("H55").Value = 1041.9166
$ws.Range("I55").Value = 954.3333
$ws.Range("J55").Value = 1129.5
$ws.Range("K55").Value = 954.3333
$ws.Range("L55").Value = 1129.5
$ws.Range("M55").Value = -781.3333
$ws.Range("N55").Value = -1475.5
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256
$ws.Range("H93").Value = 47619904
$ws.Range("I93").Value = 55556388
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 55556388
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -55555140
$ws.Range("N93").Value = -3496

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 14032.875
$ws.Range("J69").Value = 14032.875
$ws.Range("L69").Value = 14032.875
$ws.Range("N69").Value = -15530.875
$ws.Range("H72").Value = 14032.875
$ws.Range("J72").Value = 14032.875
$ws.Range("L72").Value = 42098.625
$ws.Range("N72").Value = -49586.625
$ws.Range("H100").Value = 1024.7142
$ws.Range("I100").Value = 1095
$ws.Range("K100").Value = 2190
$ws.Range("M100").Value = -1649
